$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("catalogo")

# New catalog rows describing the "semestre" / "dias no laborales" services.
$rows = @(
    @{ Row = 18; D = "crearSemestre";         E = "POST"; F = "/semestre/crearSemestre";         G = "fecha del inicio del semestre";                                   I = "JHONATAN STEVEN"; Height = 110.25 },
    @{ Row = 19; D = "listarSemestres";       E = "GET";  F = "/semestre/listarSemestres";        G = "nada";                                                            I = "JHONATAN STEVEN"; Height = 110.25 },
    @{ Row = 20; D = "eliminarSemestre";      E = "POST"; F = "/semestre/eliminarSemestre";        G = "año y periodo del semestre";                                      I = "JHONATAN STEVEN"; Height = 110.25 },
    @{ Row = 21; D = "listarDiasNoLaborales"; E = "GET";  F = "/semestre/listarDiasNoLaborales";  G = "año y periodo del semestre";                                      I = "JHONATAN STEVEN"; Height = 141.75 },
    @{ Row = 22; D = "marcarDia";             E = "POST"; F = "/semestre/marcarDia";              G = "DiaNoLaboralDTO, contiene el semestre, fecha del y causa del dia"; I = "JHONATAN STEVEN"; Height = 110.25 },
    @{ Row = 23; D = "eliminarDia";           E = "POST"; F = "/semestre/eliminarDia";            G = "id del dia no laboral";                                            I = "JHONATAN STEVEN"; Height = 110.25 }
)

# The JSON "retorno" text for each row, filled in a second pass (matches the
# order the strings were originally typed - names/paths/params first, then
# the example responses).
$responses = @(
    "{`n  `"obj`": true,`n `"mensaje`":`"semestre creo correctamente`",`n  `"codigo`": `"00`"`n} else{ `"error al crear`" }",
    "{`n  `"obj`": null,`n  `"mensaje`": `"no hay semestres`",`n  `"codigo`": `"-1`"`n} else{ la lista con lossemestres }",
    "{`n  `"obj`": true,`n  `"mensaje`": `"se elimino correctamente`",`n  `"codigo`": `"00`"`n}",
    "{`n  `"obj`": null,`n  `"mensaje`": `"no hay dias no laborales registrados`",`n  `"codigo`": `"-1`"`n} else{ la lista con los dias laborales del semestre }",
    "{`n  `"obj`": true,`n  `"mensaje`": `"el dia se marco correctamente`",`n  `"codigo`": `"00`"`n}",
    "{`n  `"obj`": true,`n  `"mensaje`": `"se elimino correctamente`",`n  `"codigo`": `"00`"`n}"
)

# First pass: Nombre servicio / Metodo / path / parametros / creado por.
foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 9).Value = $r.I
}

# Second pass: retorno (json de ejemplo), wrapped + sized like the rest of
# the table.
for ($idx = 0; $idx -lt $rows.Count; $idx++) {
    $r = $rows[$idx]
    $h = $ws.Cells.Item($r.Row, 8)
    $h.Value = $responses[$idx]
    $h.WrapText = $true
    $ws.Rows.Item($r.Row).RowHeight = $r.Height
}

# Selection / scroll position moved to show the newly added rows.
$ws.Range("H22").Select()
$excel.ActiveWindow.ScrollRow = 21

# Sheet author on this machine (reflected in the workbook's recorded save path).
$wb.Author = "JHONATAN VANEGAS"
